$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a temporary text format on price cells whose new values look like plain
# numbers, so Excel stores them as text (matching the source data) instead of
# coercing them to numeric values. The format is cleared again right after so the
# cells keep their original (default) style.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"

$ws.Range("D2").Value = "69.614.59"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "3.673.78"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "622.48"
$ws.Range("E5").Value = "  -7.37%  "
$ws.Range("D6").Value = "159.39"
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("E9").Value = "  -2.12%  "
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("D11").Value = "0.440"
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("E12").Value = "  -2.89%  "
$ws.Range("D13").Value = "4.292.41"
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("E14").Value = "  -1.69%  "
$ws.Range("D15").Value = "3.663.47"
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").Value = "69.646.04"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").Value = "6.51"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").Value = "15.86"
$ws.Range("E19").Value = "  -2.55%  "
$ws.Range("D20").Value = "10.35"
$ws.Range("E20").Value = "  +5.13%  "
$ws.Range("D21").Value = "470.00"
$ws.Range("E21").Value = "  -1.03%  "
$ws.Range("D22").Value = "0.649"
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("D23").Value = "79.65"
$ws.Range("E23").Value = "  -1.14%  "
$ws.Range("D24").Value = "3.818.75"
$ws.Range("E24").Value = "  -0.91%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  -3.50%  "
$ws.Range("D27").Value = "11.05"
$ws.Range("E27").Value = "  +0.70%  "
$ws.Range("D28").Value = "8.73"
$ws.Range("E28").Value = "  -4.42%  "
$ws.Range("E29").Value = "  -3.17%  "
$ws.Range("D30").Value = "1.67"
$ws.Range("E30").Value = "  -4.46%  "
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("E32").Value = "  -2.09%  "
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").Value = "0.163"
$ws.Range("E34").Value = "  -2.30%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "6.40"
$ws.Range("E35").Value = "  -3.45%  "
$ws.Range("D36").Value = "3.673.74"
$ws.Range("E36").Value = "  -0.64%  "
$ws.Range("E37").Value = "  -3.46%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "178.10"
$ws.Range("E39").Value = "  +3.20%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").Value = "5.78"
$ws.Range("E41").Value = "  -5.28%  "
$ws.Range("D42").Value = "2.22"
$ws.Range("E42").Value = "  -1.51%  "
$ws.Range("E43").Value = "  -2.42%  "
$ws.Range("E44").Value = "  -2.04%  "
$ws.Range("D45").Value = "46.70"
$ws.Range("D46").Value = "28.79"
$ws.Range("E46").Value = "  +3.55%  "
$ws.Range("E47").Value = "  -2.72%  "
$ws.Range("E48").Value = "  -0.60%  "
$ws.Range("E49").Value = "  -6.79%  "
$ws.Range("E50").Value = "  -5.11%  "
$ws.Range("E51").Value = "  -6.36%  "

$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D46").ClearFormats()
